# Add 6 new rows (52-57) of cumulative Covid death data to Sheet1,
# continuing the existing daily series (date serial, DeathCovid,
# DeathWithCovid, Total), and move the view/selection to follow the
# newly appended data, as in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data rows: Date (serial), DeathCovid, DeathWithCovid, Total
$newRows = @(
    @(44169, 981, 253, 1234),
    @(44170, 996, 253, 1249),
    @(44171, 1018, 256, 1274),
    @(44172, 1046, 270, 1316),
    @(44173, 1084, 274, 1358),
    @(44174, 1104, 280, 1384)
)

$startRow = 52
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$lastRow = $startRow + $newRows.Count - 1

# Update the view so it is scrolled down to show the newly added rows,
# and the new bottom-right cell is the active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D$lastRow").Select()

$wb.Save()
